# "Generate Report for Handback"
#
# The localization-status report previously showed the 6fb9166a... file as
# "Ready for handoff" everywhere that shared string was used (the Overview
# sheet's zh-cn/de-de summary columns, and the per-language "Status" column).
# A handback transform failure is now recorded instead, and the per-language
# sheets' "Error Detail" column is populated with the failure reason. The
# "Error Detail" column is also widened so the message is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# "Ready for handoff" -> "Handback transform failed" for the
# 6fb9166a-8169-4d9b-ab58-a6e4ce987ae0 row, everywhere it is shown.
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# Populate the "Error Detail" column (P) for that row with the handback
# transform failure reason, per target language.
$zhcn.Range("P3").Value = "Handback file name: zziqb0cs.idn is different with handoff file name: 6fb9166a-8169-4d9b-ab58-a6e4ce987ae0.16581e8d7490cb38fe6cd4a61f7fc513438e1774.zh-cn."
$dede.Range("P3").Value = "Handback file name: zziqb0cs.idn is different with handoff file name: 6fb9166a-8169-4d9b-ab58-a6e4ce987ae0.16581e8d7490cb38fe6cd4a61f7fc513438e1774.de-de."

# Widen the "Error Detail" column so the failure message is readable.
# ColumnWidth is in "characters" units, which Excel offsets from the
# serialized XML width by 5/6 of a character; 40 - 5/6 round-trips to an
# XML column width of exactly 40.
$zhcn.Columns.Item(16).ColumnWidth = 40 - 5/6
$dede.Columns.Item(16).ColumnWidth = 40 - 5/6
